$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (owner table header row) ---
# B2 was "owner_id" (duplicate/no-op column) -> becomes "owner_email"
$ws.Range("B2").Value = "owner_email"
# F2 was "owner_email" -> becomes "owner_rs" (previously in G2)
$ws.Range("F2").Value = "owner_rs"

# --- Row 3 (owner table sample row) ---
# B3 was "owner_id" -> becomes "owner_email"
$ws.Range("B3").Value = "owner_email"
# F3 was "owner_email" -> becomes "rs_id" (previously in G3), and must take on G3's highlighted style
$ws.Range("G3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "rs_id"

# --- Drop the now-redundant column G across rows 1-3 ---
$ws.Range("G1:G3").ClearContents()

# --- Update the active-cell selection to match the new layout ---
$ws.Range("H8").Select()
